# Apply data updates to the United Arab Emirates M2 dataset sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to existing rows (B = M2 AED, D = M2 USD recomputed with new B) ---
$updates = @(
    @{ Row = 368; B = 1486571000000;    D = 404734886126.9552 },
    @{ Row = 372; B = 1496040000000;    D = 407301833626.0056 },
    @{ Row = 379; B = 1485900000000;    D = 404555503523.0825 },
    @{ Row = 380; B = 1498300000000;    D = 407931563987.2363 },
    @{ Row = 381; B = 1517100000000;    D = 413038299494.2861 },
    @{ Row = 382; B = 1563100000000;    D = 425567829892.662  },
    @{ Row = 383; B = 1563400000000;    D = 425647189827.3752 },
    @{ Row = 384; B = 1570300000000;    D = 427524601790.632  },
    @{ Row = 385; B = 1593400000000;    D = 433813730174.6118 },
    @{ Row = 386; B = 1567300000000;    D = 426707831870.6346 },
    @{ Row = 387; B = 1568100000000;    D = 426923893689.6613 },
    @{ Row = 388; B = 1622200000000;    D = 441656525065.41   },
    @{ Row = 389; B = 1606000000000;    D = 437240592588.9625 },
    @{ Row = 390; B = 1627700000000;    D = 443153942700.6336 },
    @{ Row = 391; B = 1645600000000;    D = 448021257595.8879 },
    @{ Row = 392; B = 1629400000000;    D = 443614969214.5804 },
    @{ Row = 393; B = 1678100000000;    D = 456875111455.3654 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}

# --- New rows appended at the bottom of the table ---
$newRows = @(
    @{ Row = 400; A = 45078; B = 1855306000000; C = 0.2722551575336405; D = 505116627303.1085 },
    @{ Row = 401; A = 45108; B = 1858844000000; C = 0.2722521926510966; D = 506074354796.335  },
    @{ Row = 402; A = 45139; B = 1860330000000; C = 0.272253304474483;  D = 506480989913.015  }
)

foreach ($r in $newRows) {
    # Copy the date cell's format (from the row above) before writing the value,
    # so the new cell picks up the same style index (border/alignment/number format).
    $ws.Range("A" + ($r.Row - 1)).Copy()
    $ws.Range("A" + $r.Row).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
